# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Reorder: "Asturias" moves ahead of "Segovia"/"Leon" in the ranking table
#    (rows 22-24 shift labels while keeping the table sorted by total cases)
#    Row 22 gets new data for Asturias; rows 23 & 24 take over the previous
#    values that belonged to the old row 22 (Segovia) and row 23 (Leon).
$ws.Range("A22").Value = "Asturias"
$ws.Range("B22").Value = 2298
$ws.Range("C22").Value = 596
$ws.Range("D22").Value = 1506
$ws.Range("E22").Value = 196

$ws.Range("A23").Value = "Segovia"
$ws.Range("B23").Value = 2285
$ws.Range("C23").Value = 636
$ws.Range("D23").Value = 1469
$ws.Range("E23").Value = 167

$ws.Range("A24").Value = "Leon"
$ws.Range("B24").Value = 2285
$ws.Range("C24").Value = 1031
$ws.Range("D24").Value = 927
$ws.Range("E24").Value = 290

# 2. Row 30 (Murcia) updated case counts
$ws.Range("B30").Value = 1654
$ws.Range("C30").Value = 652
$ws.Range("D30").Value = 886
$ws.Range("E30").Value = 116

# 3. Row 54 (Mallorca) updated case counts
$ws.Range("B54").Value = 113
$ws.Range("C54").Value = 54
$ws.Range("D54").Value = 55

# 4. Row 55 (Ceuta) updated case counts
$ws.Range("B55").Value = 107
$ws.Range("C55").Value = 44
$ws.Range("D55").Value = 61

# 5. Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 22:22"
